# Update "想去人数" (column F) figures across the three data sheets
# (展览, 演出, 全部类型) to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 141
$ws1.Range("F5").Value = 215
$ws1.Range("F7").Value = 1259
$ws1.Range("F8").Value = 424
$ws1.Range("F9").Value = 207
$ws1.Range("F10").Value = 61
$ws1.Range("F12").Value = 389
$ws1.Range("F13").Value = 436
$ws1.Range("F14").Value = 808
$ws1.Range("F15").Value = 194
$ws1.Range("F16").Value = 741
$ws1.Range("F17").Value = 304
$ws1.Range("F18").Value = 87
$ws1.Range("F19").Value = 1044
$ws1.Range("F20").Value = 486
$ws1.Range("F21").Value = 286
$ws1.Range("F25").Value = 50
$ws1.Range("F27").Value = 42

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 375
$ws2.Range("F5").Value = 46
$ws2.Range("F8").Value = 86
$ws2.Range("F12").Value = 131

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 141
$ws4.Range("F7").Value = 215
$ws4.Range("F9").Value = 1259
$ws4.Range("F10").Value = 424
$ws4.Range("F11").Value = 207
$ws4.Range("F13").Value = 61
$ws4.Range("F14").Value = 375
$ws4.Range("F16").Value = 46
$ws4.Range("F17").Value = 389
$ws4.Range("F20").Value = 436
$ws4.Range("F21").Value = 808
$ws4.Range("F22").Value = 194
$ws4.Range("F23").Value = 741
$ws4.Range("F24").Value = 304
$ws4.Range("F25").Value = 87
$ws4.Range("F26").Value = 1044
$ws4.Range("F27").Value = 486
$ws4.Range("F28").Value = 86
$ws4.Range("F30").Value = 286
$ws4.Range("F36").Value = 50
$ws4.Range("F37").Value = 131
$ws4.Range("F42").Value = 42
